$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of A2:A6 (values only), keeping existing formatting/styles.
$ws.Range("A2:A6").ClearContents()

# Update the selection to match the new active cell / selected range.
$ws.Range("A2:A6").Select()
